$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = "Enter `"First Name`"    "
$ws.Range("A7").Value = $text
$ws.Range("A8").Value = $text

$ws.Range("A9").Select()
